# Apply the "想去人数" (want-to-go count) refresh and the updated venue
# address for 南昌·DSL国风动漫游戏嘉年华 to both the "展览" sheet and the
# "全部类型" sheet, which carry duplicate copies of the same table.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 3104
    $ws.Range("F5").Value = 2735
    $ws.Range("F9").Value = 1451
    $ws.Range("F11").Value = 65
    $ws.Range("F12").Value = 19
    $ws.Range("F15").Value = 372
    $ws.Range("F16").Value = 331
    $ws.Range("F17").Value = 45
    $ws.Range("F20").Value = 76
    $ws.Range("D21").Value = "九龙大道1388号（上饶街与九龙大道交叉口西北100米） 中国南昌虚拟现实VR产业基地"
    $ws.Range("F22").Value = 2668
}
